$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style from H1 to I1:J1 so new headers match formatting
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill I2:I38 and J2:J38 with data values
$iVals = @(5,6,6,7,7,9,7,6,5,6,6,4,6,9,7,7,6,8,8,6,7,8,8,8,7,9,8,6,7,8,9,9,4,4,7,7,7)
$jVals = @(5,7,6,7,7,9,7,7,5,7,6,4,6,9,7,7,6,8,8,7,7,8,8,8,7,9,8,6,7,8,9,9,4,4,7,7,7)
for ($r = 2; $r -le 38; $r++) {
    $ws.Cells.Item($r, 9).Value = $iVals[$r - 2]
    $ws.Cells.Item($r, 10).Value = $jVals[$r - 2]
}

Write-Output "done"
